$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 20 de Abril de 2020 a las 17:22"

# Estados Unidos (row 4) - updated case counts
$ws.Range("B4").Value = 767379
$ws.Range("C4").Value = 2743
$ws.Range("E4").Value = 655348
$ws.Range("G4").Value = 175
$ws.Range("H4").Value = 40750

# Row 8 - updated case counts
$ws.Range("B8").Value = 146200
$ws.Range("C8").Value = 458
$ws.Range("E8").Value = 50031
$ws.Range("G8").Value = 27
$ws.Range("H8").Value = 4669

# Countries re-ranked: Chile overtakes Arabia Saudita in total cases.
# Row 28 now shows Chile with fresh updated figures.
$ws.Range("A28").Value = "Chile"
$ws.Range("B28").Value = 10507
$ws.Range("C28").Value = 419
$ws.Range("D28").Value = 4676
$ws.Range("E28").Value = 5692
$ws.Range("F28").Value = 377
$ws.Range("H28").Value = 139

# Row 29 now shows Arabia Saudita (with the figures Arabia Saudita had before the swap).
$ws.Range("A29").Value = "Arabia Saudita"
$ws.Range("B29").Value = 10484
$ws.Range("C29").Value = 1122
$ws.Range("D29").Value = 1490
$ws.Range("E29").Value = 8891
$ws.Range("F29").Value = 88
$ws.Range("G29").Value = 6
$ws.Range("H29").Value = 103

# Row 61 - updated case counts
$ws.Range("B61").Value = 2245
$ws.Range("C61").Value = 10
$ws.Range("E61").Value = 1860
$ws.Range("F61").Value = 61
$ws.Range("G61").Value = 3
$ws.Range("H61").Value = 116

# Row 89 - updated case counts
$ws.Range("B89").Value = 772
$ws.Range("C89").Value = 5
$ws.Range("E89").Value = 679

# Row 94 - updated case counts
$ws.Range("E94").Value = 542
$ws.Range("G94").Value = 1
$ws.Range("H94").Value = 6

# Row 114 - updated case counts
$ws.Range("D114").Value = 88
$ws.Range("E114").Value = 218
